# Add an "allow_choice_duplicates" = "Yes" setting column to the XLSForm
# "settings" sheet, wrap the new header text, size the new column, and
# make "settings" the active sheet/selection (as in the upstream commit
# that adds the upload-form test fixtures).

$wb = $excel.ActiveWorkbook
$settings = $wb.Worksheets.Item("settings")

# New column C: header + value.
$settings.Range("C1").Value = "allow_choice_duplicates"
$settings.Range("C2").Value = "Yes"

# Match the widened column and wrapped header cell from the target sheet.
$settings.Columns.Item(3).ColumnWidth = 24.45
$settings.Range("C1").WrapText = $true

# "settings" becomes the active/selected sheet, with C2 selected.
$settings.Activate()
[void]$settings.Range("C2").Select()
